$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 217, shifting existing rows 217..295 down to 218..296
$ws.Rows.Item(217).Insert()

# Populate the new row 217 with the new record's data.
# Columns A,B,C,E,F,G,H,I,J,R are constant/shared with the surrounding "Durazno" rows,
# so copy them from the (now shifted) row 218 which still holds the old row 217 values.
$ws.Range("A217").Value = $ws.Range("A218").Value2
$ws.Range("B217").Value = $ws.Range("B218").Text
$ws.Range("C217").Value = $ws.Range("C218").Text
$ws.Range("D217").Value = 44876
$ws.Range("E217").Value = $ws.Range("E218").Value2
$ws.Range("F217").Value = $ws.Range("F218").Text
$ws.Range("G217").Value = $ws.Range("G218").Value2
$ws.Range("H217").Value = $ws.Range("H218").Text
$ws.Range("I217").Value = $ws.Range("I218").Value2
$ws.Range("J217").Value = $ws.Range("J218").Text
$ws.Range("K217").Value = "Early Majestic"
$ws.Range("L217").Value = "Primera"
$ws.Range("M217").Value = 160
$ws.Range("N217").Value = 20000
$ws.Range("O217").Value = 22000
$ws.Range("P217").Value = 21000
$ws.Range("Q217").Value = "`$/bandeja 10 kilos granel"
$ws.Range("R217").Value = $ws.Range("R218").Text
$ws.Range("S217").Value = 2100
$ws.Range("T217").Value = 10
